# Edit sheet Card24 by admin
# Remove column O ("Servised by") entirely - header text and the empty
# placeholder cells below it - from the Card24 worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$ws.Columns.Item(15).Delete()
